$wb = $excel.ActiveWorkbook

# --- Summary sheet: widen the selected/used range in the sheet view ---
$wsSummary = $wb.Worksheets.Item("Summary")
$null = $wsSummary.Activate()
$null = $wsSummary.Range("A7:XFD14").Select()

# --- Repayment schedule sheet: add column O (copy format from column N) ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Rows 2-3: formats only (values stay blank, matching column N in those rows)
for ($r = 2; $r -le 3; $r++) {
    $nCell = $wsSchedule.Cells.Item($r, 14)
    $oCell = $wsSchedule.Cells.Item($r, 15)
    $null = $nCell.Copy()
    $null = $oCell.PasteSpecial(-4122)
}

# Rows 4-14: copy format from column N and set value to 0
for ($r = 4; $r -le 14; $r++) {
    $nCell = $wsSchedule.Cells.Item($r, 14)
    $oCell = $wsSchedule.Cells.Item($r, 15)
    $null = $nCell.Copy()
    $null = $oCell.PasteSpecial(-4122)
    $oCell.Value = 0
}

$excel.CutCopyMode = $false

# --- Transactions sheet: update the transaction IDs and the active selection ---
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("A2").Value = 100
$wsTransactions.Range("A3").Value = 98
$wsTransactions.Range("A4").Value = 96

$null = $wsTransactions.Activate()
$null = $wsTransactions.Range("D4").Select()
